# SENG637-A4.pptx — slide 5 ("Assignment 4 – Part 2"), content placeholder:
#   - 1st bullet "Install Selenium and SikuliX" -> "Install Selenium"
#     (drop the " and SikuliX" run + stray endParaRPr)
#   - last bullet "Compare Selenium with SikuliX" paragraph is removed entirely
#
# The middle two (indented) bullets are left completely untouched.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(5)

# Locate the body/content placeholder by name so this isn't dependent on a
# brittle shape index.
$shape = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $candidate = $s.Shapes.Item($i)
    if ($candidate.Name -eq "Content Placeholder 2") {
        $shape = $candidate
    }
}

$tr = $shape.TextFrame.TextRange

# Rewriting the whole run collapses the multi-run first paragraph down to a
# single run and drops the no-longer-needed paragraph marks/endParaRPr, while
# still recreating each paragraph with the same "en-CA"/dirty="0" run
# formatting the deck already used everywhere else on this slide.
$tr.Text = "Install Selenium`rFollow instructions to `rDesign test cases for at least 2 functionalities (per student) of the selected websites`rAutomate your designed test cases using Selenium and add verification points to your scripts"

# Restore the indent level (PowerPoint's 1-based IndentLevel == OOXML's
# 0-based <a:pPr lvl="1"/>) on the two bullets that were already indented.
$tr.Paragraphs(3, 1).IndentLevel = 2
$tr.Paragraphs(4, 1).IndentLevel = 2
